$wb = $excel.ActiveWorkbook

# --- On the "search_product" sheet, change the searched product from
#     "samsung TV" to "65-inch TV" ---
$wsSearch = $wb.Worksheets.Item("search_product")
$wsSearch.Range("A2").Value = "65-inch TV"

# --- Make "search_product" the active sheet/tab (previously "sign_in" was
#     the selected tab) ---
$wsSearch.Activate()
